$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.474.65'
$ws.Range("E2").Value = '  -0.84%  '

$ws.Range("D3").Value = '1.875.85'
$ws.Range("E3").Value = '  -0.75%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4713'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.69%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2872'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06512'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.91'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '100.82'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.92%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07806'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.10%  '

$ws.Range("D13").Value = '1.874.10'
$ws.Range("E13").Value = '  -0.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7281'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.175'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.41%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '284.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("D17").Value = '30.475.42'
$ws.Range("E17").Value = '  -0.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9992'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007488'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.64%  '

$ws.Range("D21").Value = '2.117.78'
$ws.Range("E21").Value = '  -1.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.339'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9996'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.348'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.057'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.900'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.09681'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.321'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.498'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.229'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.69%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.150'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04810'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.11%  '

$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.766'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.68%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6906'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01903'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.848'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.85'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.29%  '

$ws.Range("E41").Value = '  -1.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.962'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.82%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4222'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9987'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8248'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.737'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.016'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05758'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '883.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.81%  '
